# Add a new input row for district heating demand (dh_demand)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67
$ws.Cells.Item($row, 1).Value = "district heating demand "
$ws.Cells.Item($row, 2).Value = "input"
$ws.Cells.Item($row, 3).Value = "dh_demand"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = "MWh"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = [double]"9.9999999999999998E+101"

# Match the author's final selection/viewport on the newly added cell
$ws.Range("H67").Select() | Out-Null
